# Update the "想去人数" (want-to-go count) figures for several events.
# These values are duplicated across the "展览" and "全部类型" worksheets,
# since "全部类型" aggregates all event rows shown on the other sheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row -> new value, per sheet (row numbers taken directly from each sheet)
$exhibitUpdates = @{
    7  = 5580
    10 = 3819
    11 = 68
    12 = 20
    20 = 552
    23 = 5242
    28 = 7777
    35 = 1185
    41 = 14
}

$allUpdates = @{
    9  = 5580
    10 = 3819
    11 = 68
    12 = 20
    20 = 552
    24 = 5242
    29 = 7777
    36 = 1185
    39 = 14
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
